function Merge-Text($tr, $oldText, $newText) {
    $full = $tr.Text
    $idx = $full.IndexOf($oldText)
    if ($idx -lt 0) {
        Write-Host "WARN: substring not found -> $oldText"
        return
    }
    $len = $oldText.Length
    $sub = $tr.Characters($idx + 1, $len)
    $sub.Text = $newText
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# Slide 1: Install MOOSE (liveCD/USB instructions)
# ---------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$tr1 = $s1.Shapes.Item(1).TextFrame.TextRange

Merge-Text $tr1 ") on Ubuntu 12.04. Boot off it, on your machine, or in a virtual machine (virtualbox.org)." ") on Ubuntu 12.04. Boot off it, on your machine, or using virtualbox.org. User:moose, passwd:moose."

Merge-Text $tr1 "Quick start MOOSE:" "Quick start MOOSE:"

Merge-Text $tr1 "You can ``stop' a simulation in-between and ``run' again for another ``Run Time'. ``Reset' and ``run' to start from t=0." "You can ‘stop’ a simulation in-between and ‘run’ again for another ‘Run Time’. ‘Reset’ and ‘run’ to start from t=0."

# ---------------------------------------------------------------
# Slide 2: Changing Parameters
# ---------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$tr2 = $s2.Shapes.Item(1).TextFrame.TextRange

Merge-Text $tr2 "Changing Parameters:" "Changing Parameters:"

Merge-Text $tr2 "An empty unlinked box near a filled ellipse signifies an explicit enzyme complex for an explicit enzyme reaction." "An empty unlinked box near a filled ellipse signifies an explicit enzyme complex for an explicit enzyme reaction."

# ---------------------------------------------------------------
# Slide 3: Plotting / Switching Solvers / Neuron-Electrical headers
# ---------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$tr3 = $s3.Shapes.Item(1).TextFrame.TextRange

Merge-Text $tr3 "Plotting:" "Plotting:"

Merge-Text $tr3 "Switching Solvers:" "Switching Solvers:"

Merge-Text $tr3 "Neuron / Electrical:" "Neuron / Electrical:"

